# Update the Maltaspor roster sheet: refresh player data and reorder rows
# to match the newly uploaded data set.
# Row 1 (header), row 3 (LaMelo Ball), row 17 and row 18 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (Player, Position, Team)
$rows = @{
    2  = @("T.J. McConnell",    "PG",       "Indiana Pacers")
    4  = @("Derrick White",     "PG,SG",    "Boston Celtics")
    5  = @("Cam Thomas",        "SG,SF",    "Brooklyn Nets")
    6  = @("Cameron Johnson",   "SF,PF",    "Brooklyn Nets")
    7  = @("Brandon Ingram",    "SG,SF,PF", "New Orleans Pelicans")
    8  = @("Cody Martin",       "SG,SF",    "Charlotte Hornets")
    9  = @("Brandon Miller",    "SG,SF",    "Charlotte Hornets")
    10 = @("Isaiah Hartenstein","C",        "Oklahoma City Thunder")
    11 = @("Damian Lillard",    "PG",       "Milwaukee Bucks")
    12 = @("Julius Randle",     "PF",       "Minnesota Timberwolves")
    13 = @("Bam Adebayo",       "C",        "Miami Heat")
    14 = @("Jared McCain",      "PG,SG",    "Philadelphia 76ers")
    15 = @("Bilal Coulibaly",   "SG,SF",    "Washington Wizards")
    16 = @("Anthony Davis",     "PF,C",     "Los Angeles Lakers")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
